$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(3)
$tr = $shape.TextFrame.TextRange
$full = $tr.Characters(1, $tr.Length)
$full.Text = "Followed by a picture"
